# This script merges the two separate headings "5) Add Tracks to a
# Playlist" and "Remove Tracks from a Playlist" into a single heading
# ("5) Add Tracks to a Playlist Remove Tracks from a Playlist") that now
# occupies the previously-empty heading paragraph right before them,
# wraps the new title text in a bookmark (matching the document's
# existing pattern for other numbered headings), empties out the old
# "5) Add Tracks to a Playlist" paragraph, and removes the old
# "Remove Tracks from a Playlist" paragraph entirely.

$d = $word.ActiveDocument

# Locate the three paragraphs involved by searching for recognizable
# text, since paragraph indices could in principle shift.
$emptyHeadingIndex = 0
$addParaIndex = 0
$removeParaIndex = 0

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($t -eq "5) Add Tracks to a Playlist") {
        $addParaIndex = $i
    } elseif ($t -eq "Remove Tracks from a Playlist") {
        $removeParaIndex = $i
    }
}

if ($addParaIndex -eq 0 -or $removeParaIndex -eq 0) {
    Write-Host "ERROR: could not locate target paragraphs"
}

$emptyHeadingIndex = $addParaIndex - 1

# --- Step 1: fill the empty heading paragraph with the merged title ---
$emptyPara = $d.Paragraphs.Item($emptyHeadingIndex)
$insertionPoint = $d.Range($emptyPara.Range.Start, $emptyPara.Range.Start)

$numberPrefix = "5) "
$mergedTitle = "Add Tracks to a Playlist Remove Tracks from a Playlist"
$fullText = $numberPrefix + $mergedTitle

# Insert the whole string in one shot (avoids the engine carrying the
# "preserve whitespace" flag from the "5) " run onto the rest of the
# text when runs are later split apart).
$insertionPoint.InsertAfter($fullText)

# Re-fetch the paragraph and wrap the title portion (after "5) ") in a
# bookmark; this naturally splits it into its own run, matching the
# style used for the sibling "3)" / "4)" headings elsewhere in the
# document.
$filledPara = $d.Paragraphs.Item($emptyHeadingIndex)
$paraStart = $filledPara.Range.Start
$bookmarkRange = $d.Range($paraStart + $numberPrefix.Length, $paraStart + $fullText.Length)
$d.Bookmarks.Add("__DdeLink__107_577582903", $bookmarkRange) | Out-Null

# --- Step 2: clear the old "5) Add Tracks to a Playlist" paragraph ---
$oldAddPara = $d.Paragraphs.Item($addParaIndex)
$oldAddPara.Range.Text = ""

# --- Step 3: delete the old "Remove Tracks from a Playlist" paragraph entirely ---
$oldRemovePara = $d.Paragraphs.Item($removeParaIndex)
$oldRemovePara.Range.Delete()

# --- Accompanying style-sheet touch-ups that came along with this edit ---

# Normal style: stop allowing punctuation to overflow the text margins.
$normalStyle = $d.Styles.Item("Normal")
$normalStyle.ParagraphFormat.HangingPunctuation = 0

# New (unused) character "ListLabel" styles that the authoring tool
# added to the stylesheet alongside the edit.
foreach ($n in 40..48) {
    $s = $d.Styles.Add("ListLabel $n", 2)
    $s.QuickStyle = $true
    $s.Font.NameBi = "OpenSymbol"
}

$s49 = $d.Styles.Add("ListLabel 49", 2)
$s49.QuickStyle = $true
$s49.Font.Bold = $false
$s49.Font.BoldBi = $false

$s50 = $d.Styles.Add("ListLabel 50", 2)
$s50.QuickStyle = $true
$s50.Font.Bold = $true
$s50.Font.BoldBi = $true
$s50.Font.Size = 12
$s50.Font.SizeBi = 12

$s51 = $d.Styles.Add("ListLabel 51", 2)
$s51.QuickStyle = $true

Write-Host "Edit complete"
